# Update Name of Algo - update imputed values in column C for several rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -12.0398
    7   = -13.28170000000001
    16  = -14.30599999999999
    28  = -13.4965
    29  = -11.7197
    32  = -12.3679
    40  = -12.5424
    52  = -11.1486
    57  = -14.05
    66  = -11.3896
    100 = -11.6869
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
